# Applies the "Issues fixes and reports" commit:
#  - AMSIN sheet: append rows 58 and 59
#  - BETA  sheet: append row 29
#  - AMS   sheet: append row 28 (and normalize row 27's value/formatting)
#
# Column layout for every sheet touched here:
#   A Run Date (text, looks like a date -> must be written as TEXT, not an
#     auto-converted date serial number)
#   B Run Time (numeric date/time serial)
#   C Sprint Name (text)
#   D Total Cases (numeric)
#   E Pass Cases (numeric)
#   F Fail Cases (numeric)
#   G Time Taken (numeric)

$wb = $excel.ActiveWorkbook

# A scratch cell used to safely get a "date looking" string into the grid as
# plain text. Writing a formula that RETURNS a string (instead of assigning a
# literal string straight to the destination cell) avoids Excel's "smart"
# date auto-detection, and avoids permanently creating a brand new number
# format / style entry the way forcing NumberFormat = "@" would.
function Set-TextDateCell($ws, $row, $col, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-TextCell($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}

# Same as the plain setters above, but clears the existing cell content
# first. Used when touching a cell that already exists in the sheet (as
# opposed to a cell in a brand new row) so the saved style index is written
# out explicitly instead of being collapsed into "inherit the column
# default".
function Set-NumCellExplicit($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).ClearContents()
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-TextCellExplicit($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).ClearContents()
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-TextDateCellExplicit($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).ClearContents()
    Set-TextDateCell $ws $row $col $text
}

# Run Time (column B) uses a dedicated date/time number format (style 11)
# that differs from the plain column default (style 5). Pull that exact
# formatting down from the row above before writing the new value so the
# saved style index lines up with the rest of the column.
function Set-RunTimeCell($ws, $row, $col, $val) {
    $above = $ws.Cells.Item($row - 1, $col)
    $above.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($row, $col).Value = $val
}

function Add-HistoryRow($ws, $row, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    Set-TextDateCell $ws $row 1 $runDate
    Set-RunTimeCell  $ws $row 2 $runTime
    Set-TextCell     $ws $row 3 $sprintName
    Set-NumCell      $ws $row 4 $total
    Set-NumCell      $ws $row 5 $pass
    Set-NumCell      $ws $row 6 $fail
    Set-NumCell      $ws $row 7 $timeTaken
}

# ---- AMSIN: dimension A1:G57 -> A1:G59 ----
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Add-HistoryRow $wsAmsin 58 "2023-02-17" 44974.43926104167 "173cyclefst" 105 105 0 2.74
Add-HistoryRow $wsAmsin 59 "2023-02-20" 44977.41202688657 "173fnlrun"   105 105 0 2.8

# ---- BETA: dimension A1:G28 -> A1:G29 ----
$wsBeta = $wb.Worksheets.Item("BETA")
Add-HistoryRow $wsBeta 29 "2023-02-20" 44977.59907541666 "173beta" 105 103 2 3.41

# ---- AMS: dimension A1:G27 -> A1:G28 ----
$wsAms = $wb.Worksheets.Item("AMS")

# Normalize the pre-existing row 27 (it previously lacked explicit styling on
# most cells and had a slightly different B27 raw value).
Set-TextDateCellExplicit $wsAms 27 1 "2023-01-20"
Set-NumCell              $wsAms 27 2 44946.898035625
Set-TextCellExplicit     $wsAms 27 3 "172live"
Set-NumCellExplicit      $wsAms 27 4 105
Set-NumCellExplicit      $wsAms 27 5 105
Set-NumCellExplicit      $wsAms 27 6 0
Set-NumCellExplicit      $wsAms 27 7 2.85

Add-HistoryRow $wsAms 28 "2023-02-20" 44977.83812617893 "live173" 105 105 0 2.78
